$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1716.921
$ws.Range("I17").Value = 437.5
$ws.Range("J17").Value = 1867.4412
$ws.Range("K17").Value = 1312.5
$ws.Range("L17").Value = 5602.3236
$ws.Range("M17").Value = -1144.5
$ws.Range("N17").Value = -5938.3236
$ws.Range("H62").Value = 5104.1665
$ws.Range("I62").Value = 3720.6924
$ws.Range("K62").Value = 3720.6924
$ws.Range("M62").Value = -3096.6924
$ws.Range("H65").Value = 5104.1665
$ws.Range("I65").Value = 3720.6924
$ws.Range("K65").Value = 18603.462
$ws.Range("M65").Value = -15483.462
$ws.Range("H70").Value = 12503677
$ws.Range("J70").Value = 20005364
$ws.Range("L70").Value = 60016092
$ws.Range("N70").Value = -60016632
$ws.Range("H73").Value = 12503677
$ws.Range("J73").Value = 20005364
$ws.Range("L73").Value = 60016092
$ws.Range("N73").Value = -60017964
$ws.Range("H112").Value = 1897.8889
$ws.Range("I112").Value = 865.8
$ws.Range("K112").Value = 2597.4
$ws.Range("M112").Value = -1489.4
$ws.Range("H113").Value = 6833.7334
$ws.Range("I113").Value = 2600
$ws.Range("J113").Value = 7892.1665
$ws.Range("K113").Value = 2600
$ws.Range("L113").Value = 7892.1665
$ws.Range("M113").Value = 654
$ws.Range("N113").Value = -14400.1665
$ws.Range("H135").Value = 971.26086
$ws.Range("I135").Value = 778
$ws.Range("K135").Value = 7002
$ws.Range("M135").Value = -4467
$ws.Range("H137").Value = 2601.1904
$ws.Range("I137").Value = 1809.3158
$ws.Range("K137").Value = 5427.9474
$ws.Range("M137").Value = -2877.9474
$ws.Range("H138").Value = 3311.1765
$ws.Range("J138").Value = 3991.3333
$ws.Range("L138").Value = 11973.9999
$ws.Range("N138").Value = -22253.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1739.4286
$ws.Range("I32").Value = 1341.6061
$ws.Range("K32").Value = 1341.6061
$ws.Range("M32").Value = -1054.6061
$ws.Range("H45").Value = 58826412
$ws.Range("I45").Value = 83334680
$ws.Range("K45").Value = 83334680
$ws.Range("M45").Value = -83334303
$ws.Range("H61").Value = 8706.764999999999
$ws.Range("I61").Value = 6460.926
$ws.Range("K61").Value = 6460.926
$ws.Range("M61").Value = -6248.926
$ws.Range("H102").Value = 1806
$ws.Range("I102").Value = 1828.7858
$ws.Range("K102").Value = 1828.7858
$ws.Range("M102").Value = -206.7858000000001
$ws.Range("H122").Value = 2421.8
$ws.Range("I122").Value = 1136.6666
$ws.Range("K122").Value = 3409.9998
$ws.Range("M122").Value = -959.9998000000001
$ws.Range("H136").Value = 8706.764999999999
$ws.Range("I136").Value = 6460.926
$ws.Range("K136").Value = 19382.778
$ws.Range("M136").Value = -16832.778

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1226.2
$ws.Range("I94").Value = 1164.1111
$ws.Range("J94").Value = 1319.3334
$ws.Range("K94").Value = 1164.1111
$ws.Range("L94").Value = 1319.3334
$ws.Range("M94").Value = -713.1111000000001
$ws.Range("N94").Value = -2221.3334
$ws.Range("H105").Value = 18571.285
$ws.Range("I105").Value = 13000
$ws.Range("K105").Value = 13000
$ws.Range("M105").Value = -11253
$ws.Range("H107").Value = 1401.2307
$ws.Range("I107").Value = 542.4286
$ws.Range("K107").Value = 542.4286
$ws.Range("M107").Value = 1377.5714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 27662.953
$ws.Range("J31").Value = 110214.5
$ws.Range("L31").Value = 110214.5
$ws.Range("N31").Value = -110804.5
$ws.Range("H34").Value = 27662.953
$ws.Range("J34").Value = 110214.5
$ws.Range("L34").Value = 110214.5
$ws.Range("N34").Value = -110618.5
$ws.Range("H74").Value = 266665.66
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 266665.66
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 266665.66
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -268413.66
$ws.Range("H77").Value = 266665.66
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 266665.66
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 799996.98
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -808732.98
$ws.Range("H86").Value = 12266.733
$ws.Range("I86").Value = 9499.25
$ws.Range("K86").Value = 9499.25
$ws.Range("M86").Value = -8376.25
$ws.Range("H89").Value = 12266.733
$ws.Range("I89").Value = 9499.25
$ws.Range("K89").Value = 47496.25
$ws.Range("M89").Value = -41880.25
$ws.Range("H122").Value = 8729.733
$ws.Range("J122").Value = 14056
$ws.Range("L122").Value = 42168
$ws.Range("N122").Value = -47068
$ws.Range("H141").Value = 168738.25
$ws.Range("J141").Value = 189986.58
$ws.Range("L141").Value = 189986.58
$ws.Range("N141").Value = -200346.58

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 72.35714
$ws.Range("I33").Value = 75.888885
$ws.Range("J33").Value = 66
$ws.Range("K33").Value = 455.33331
$ws.Range("L33").Value = 396
$ws.Range("M33").Value = -172.33331
$ws.Range("N33").Value = -962
$ws.Range("H42").Value = 11334.667
$ws.Range("J42").Value = 14502
$ws.Range("L42").Value = 43506
$ws.Range("N42").Value = -44574
$ws.Range("H69").Value = 8336.888999999999
$ws.Range("I69").Value = 5007.6665
$ws.Range("K69").Value = 15022.9995
$ws.Range("M69").Value = -14211.9995
$ws.Range("H72").Value = 8336.888999999999
$ws.Range("I72").Value = 5007.6665
$ws.Range("K72").Value = 45068.9985
$ws.Range("M72").Value = -41012.9985
$ws.Range("H121").Value = 2436.3684
$ws.Range("I121").Value = 1592.5
$ws.Range("J121").Value = 2825.8462
$ws.Range("K121").Value = 4777.5
$ws.Range("L121").Value = 8477.5386
$ws.Range("M121").Value = -3467.5
$ws.Range("N121").Value = -11097.5386
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("H138").Value = 4010.875
$ws.Range("I138").Value = 1012.4286
$ws.Range("K138").Value = 3037.2858
$ws.Range("M138").Value = 2102.7142
$ws.Range("H139").Value = 6957.75
$ws.Range("I139").Value = 1632.5
$ws.Range("J139").Value = 12283
$ws.Range("K139").Value = 4897.5
$ws.Range("L139").Value = 36849
$ws.Range("M139").Value = 242.5
$ws.Range("N139").Value = -47129

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5474.846
$ws.Range("I122").Value = 5018.3335
$ws.Range("J122").Value = 6502
$ws.Range("K122").Value = 15055.0005
$ws.Range("L122").Value = 19506
$ws.Range("M122").Value = -12605.0005
$ws.Range("N122").Value = -24406

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 5071.4287
$ws.Range("I68").Value = 2123
$ws.Range("K68").Value = 2123
$ws.Range("M68").Value = -1374
$ws.Range("H71").Value = 5071.4287
$ws.Range("I71").Value = 2123
$ws.Range("K71").Value = 10615
$ws.Range("M71").Value = -6871
$ws.Range("H93").Value = 6067.4287
$ws.Range("I93").Value = 4619.5
$ws.Range("J93").Value = 7998
$ws.Range("K93").Value = 4619.5
$ws.Range("L93").Value = 7998
$ws.Range("M93").Value = -3371.5
$ws.Range("N93").Value = -10494
$ws.Range("H100").Value = 6974.35
$ws.Range("I100").Value = 6191.933
$ws.Range("J100").Value = 9321.6
$ws.Range("K100").Value = 6191.933
$ws.Range("L100").Value = 9321.6
$ws.Range("M100").Value = -5650.933
$ws.Range("N100").Value = -10403.6
$ws.Range("H122").Value = 5470.75
$ws.Range("J122").Value = 4963.3335
$ws.Range("L122").Value = 14890.0005
$ws.Range("N122").Value = -19790.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5677.8335
$ws.Range("J81").Value = 10501
$ws.Range("L81").Value = 21002
$ws.Range("N81").Value = -23124
$ws.Range("H84").Value = 5677.8335
$ws.Range("J84").Value = 10501
$ws.Range("L84").Value = 105010
$ws.Range("N84").Value = -115618

